$wb = $excel.ActiveWorkbook

# --- Carabao_cup sheet: "Swansea City" -> "Swansea" for the 12/13 season ---
$ws2 = $wb.Worksheets.Item("Carabao_cup")
$ws2.Range("B12").Value = "Swansea"

# --- Fa_cup sheet: "Wigan Athletic" -> "Wigan" for the 12/13 season ---
$ws1 = $wb.Worksheets.Item("Fa_cup")
$ws1.Range("B11").Value = "Wigan"

# Row 32 on Fa_cup loses its extra (applyNumberFormat) cell style
$ws1.Range("A32:B32").ClearFormats()

# --- Update the saved selection / scroll state on each sheet ---
$ws1.Activate()
$ws1.Range("B12").Select() | Out-Null

$ws2.Activate()
$ws2.Range("B13").Select() | Out-Null

# Restore Fa_cup as the active / selected tab
$ws1.Activate()
